$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append " *" to each of the header cells to indicate required fields
$ws.Range("A1").Value = "Fund *"
$ws.Range("B1").Value = "Name *"
$ws.Range("C1").Value = "Percentage Called *"
$ws.Range("D1").Value = "Due Date *"

# Move the active selection to D2
$ws.Range("D2").Select()
